$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 3, 5, 7 as part of re-pulled data / mean calculation refresh
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = 9
